$d = $word.ActiveDocument

# The new log entry goes right after the last paragraph in the body
# ("Ok I added thumbnail for both pictures and videos, ...").
$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count)
$last.Range.InsertParagraphAfter()

$newLast = $d.Paragraphs($d.Paragraphs.Count)
$newLast.Range.Text = "Reading .json is easy, the problem is how to read all and filter etc. Ok so to start from basic, im going to implement basic search/query first."
